$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "#6"
$ws.Range("B7").Value = "OPD service link"
$ws.Range("C7").Value = "Test if OPD service link works"
$ws.Range("D7").Value = "It will link to a page About OPD"

$ws.Range("D10").Select()
